$d = $word.ActiveDocument

# Locate the sentence "...switch the display to the video." inside the
# "6. Video Controls" section and append a new sentence describing the
# keyboard shortcuts, broken into the same separate runs as the authored
# edit (e.g. "left" / "right" isolated in their own runs).
$rng = $d.Content
$found = $rng.Find.Execute("switch the display to the video.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found text so inserts land right after it.
    $rng.Collapse(0)

    $rng.InsertAfter(" The space, ")
    $rng.Collapse(0)

    $rng.InsertAfter("left")
    $rng.Collapse(0)

    $rng.InsertAfter(" arrow, and ")
    $rng.Collapse(0)

    $rng.InsertAfter("right")
    $rng.Collapse(0)

    $rng.InsertAfter(" arrow keys also control video playback.")
    $rng.Collapse(0)
}
